$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61; this pushes the existing rows 61-204
# down to 62-205 (so the former last row duplicates into the new last row
# 205), growing the used range from A1:R204 to A1:R205.
$ws.Rows("61").Insert()

# Populate the newly-inserted row 61 with the new price-report record.
$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 44544
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = 100112037
$ws.Range("G61").Value = "Cebollín"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 180
$ws.Range("K61").Value = 6000
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = 6000
$ws.Range("N61").Value = "`$/paquete 36 unidades"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 167
$ws.Range("Q61").Value = 36
$ws.Range("R61").Value = "Hortaliza"
